# Add a "Phone" column (E) with phone numbers for each contact row,
# matching the formatting already used for the blank trailer row (A4:D4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the Phone column
$ws.Range("E1").Value = "Phone"

# Phone numbers for the two data rows
$ws.Range("E2").Value = 919167510548
$ws.Range("E3").Value = 919167510548

# Copy the existing "plain" formatting (Arial / theme color, no border) from
# the blank A4 cell onto the new Phone cells so styles are reused rather than
# duplicated.
$ws.Range("A4").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("E3").PasteSpecial(-4122)

$excel.CutCopyMode = 0
